$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 210.33333
$ws.Range("I55").Value = 147.55556
$ws.Range("J55").Value = 398.66666
$ws.Range("K55").Value = 147.55556
$ws.Range("L55").Value = 398.66666
$ws.Range("M55").Value = 66.44443999999999
$ws.Range("N55").Value = -826.66666
$ws.Range("H70").Value = 3744.1562
$ws.Range("I70").Value = 2362.9048
$ws.Range("K70").Value = 7088.714399999999
$ws.Range("M70").Value = -6818.714399999999
$ws.Range("H73").Value = 3744.1562
$ws.Range("I73").Value = 2362.9048
$ws.Range("K73").Value = 7088.714399999999
$ws.Range("M73").Value = -6152.714399999999
$ws.Range("H87").Value = 52498.5
$ws.Range("J87").Value = 52498.5
$ws.Range("L87").Value = 52498.5
$ws.Range("N87").Value = -54994.5
$ws.Range("H90").Value = 52498.5
$ws.Range("J90").Value = 52498.5
$ws.Range("L90").Value = 157495.5
$ws.Range("N90").Value = -169975.5
$ws.Range("H113").Value = 5587.6875
$ws.Range("I113").Value = 6165.3335
$ws.Range("K113").Value = 6165.3335
$ws.Range("M113").Value = -2911.3335
$ws.Range("H125").Value = 7180.3125
$ws.Range("I125").Value = 3674.5
$ws.Range("J125").Value = 10686.125
$ws.Range("K125").Value = 33070.5
$ws.Range("L125").Value = 96175.125
$ws.Range("M125").Value = -30610.5
$ws.Range("N125").Value = -101095.125
$ws.Range("H137").Value = 8033
$ws.Range("I137").Value = 12275.3
$ws.Range("K137").Value = 36825.89999999999
$ws.Range("M137").Value = -34275.89999999999
$ws.Range("H138").Value = 5437.47
$ws.Range("J138").Value = 5446.694
$ws.Range("L138").Value = 16340.082
$ws.Range("N138").Value = -26620.082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23482.203
$ws.Range("I32").Value = 19569.25
$ws.Range("J32").Value = 49568.555
$ws.Range("K32").Value = 19569.25
$ws.Range("L32").Value = 49568.555
$ws.Range("M32").Value = -19282.25
$ws.Range("N32").Value = -50142.555
$ws.Range("H45").Value = 2152.2632
$ws.Range("I45").Value = 1219
$ws.Range("J45").Value = 2485.5715
$ws.Range("K45").Value = 1219
$ws.Range("L45").Value = 2485.5715
$ws.Range("M45").Value = -842
$ws.Range("N45").Value = -3239.5715
$ws.Range("H61").Value = 4282.5
$ws.Range("I61").Value = 2510.476
$ws.Range("J61").Value = 9598.571
$ws.Range("K61").Value = 2510.476
$ws.Range("L61").Value = 9598.571
$ws.Range("M61").Value = -2298.476
$ws.Range("N61").Value = -10022.571
$ws.Range("H122").Value = 9363.8125
$ws.Range("I122").Value = 4954.125
$ws.Range("J122").Value = 13773.5
$ws.Range("K122").Value = 14862.375
$ws.Range("L122").Value = 41320.5
$ws.Range("M122").Value = -12412.375
$ws.Range("N122").Value = -46220.5
$ws.Range("H132").Value = 4247.3794
$ws.Range("I132").Value = 4149.1577
$ws.Range("K132").Value = 12447.4731
$ws.Range("M132").Value = -9917.473099999999
$ws.Range("H136").Value = 4282.5
$ws.Range("I136").Value = 2510.476
$ws.Range("J136").Value = 9598.571
$ws.Range("K136").Value = 7531.428
$ws.Range("L136").Value = 28795.713
$ws.Range("M136").Value = -4981.428
$ws.Range("N136").Value = -33895.713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1410.0714
$ws.Range("I94").Value = 1254.2
$ws.Range("K94").Value = 1254.2
$ws.Range("M94").Value = -803.2
$ws.Range("H99").Value = 7911.727
$ws.Range("I99").Value = 9419.214
$ws.Range("J99").Value = 5273.625
$ws.Range("K99").Value = 9419.214
$ws.Range("L99").Value = 5273.625
$ws.Range("M99").Value = -7921.214
$ws.Range("N99").Value = -8269.625
$ws.Range("H105").Value = 2990
$ws.Range("I105").Value = 2990
$ws.Range("K105").Value = 2990
$ws.Range("M105").Value = -1243
$ws.Range("H134").Value = 3664.6
$ws.Range("I134").Value = 4818.476
$ws.Range("K134").Value = 14455.428
$ws.Range("M134").Value = -11920.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30727822
$ws.Range("I31").Value = 8774851
$ws.Range("J31").Value = 52680796
$ws.Range("K31").Value = 8774851
$ws.Range("L31").Value = 52680796
$ws.Range("M31").Value = -8774556
$ws.Range("N31").Value = -52681386
$ws.Range("H34").Value = 30727822
$ws.Range("I34").Value = 8774851
$ws.Range("J34").Value = 52680796
$ws.Range("K34").Value = 8774851
$ws.Range("L34").Value = 52680796
$ws.Range("M34").Value = -8774649
$ws.Range("N34").Value = -52681200
$ws.Range("H105").Value = 4773.0312
$ws.Range("I105").Value = 5011.7407
$ws.Range("J105").Value = 3484
$ws.Range("K105").Value = 5011.7407
$ws.Range("L105").Value = 3484
$ws.Range("M105").Value = -3264.7407
$ws.Range("N105").Value = -6978
$ws.Range("H132").Value = 6541.3335
$ws.Range("I132").Value = 9475
$ws.Range("K132").Value = 28425
$ws.Range("M132").Value = -25895
$ws.Range("H141").Value = 135681.52
$ws.Range("J141").Value = 135681.52
$ws.Range("L141").Value = 135681.52
$ws.Range("N141").Value = -146041.52

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3248
$ws.Range("I3").Value = 3248
$ws.Range("K3").Value = 9744
$ws.Range("M3").Value = -9632
$ws.Range("H23").Value = 154.71428
$ws.Range("I23").Value = 249.33333
$ws.Range("J23").Value = 83.75
$ws.Range("K23").Value = 747.99999
$ws.Range("L23").Value = 251.25
$ws.Range("M23").Value = -512.99999
$ws.Range("N23").Value = -721.25
$ws.Range("H32").Value = 309935.1
$ws.Range("I32").Value = 755550
$ws.Range("J32").Value = 12858.5
$ws.Range("K32").Value = 2266650
$ws.Range("L32").Value = 38575.5
$ws.Range("M32").Value = -2266367
$ws.Range("N32").Value = -39141.5
$ws.Range("H50").Value = 999
$ws.Range("I50").Value = 999
$ws.Range("K50").Value = 2997
$ws.Range("M50").Value = -2516
$ws.Range("H53").Value = 999
$ws.Range("I53").Value = 999
$ws.Range("K53").Value = 2997
$ws.Range("M53").Value = -2516
$ws.Range("H122").Value = 1810
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 2087.5
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 18787.5
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -23687.5
$ws.Range("H133").Value = 4648.3
$ws.Range("J133").Value = 5357.143
$ws.Range("L133").Value = 16071.429
$ws.Range("N133").Value = -26191.429
$ws.Range("H134").Value = 4565.1333
$ws.Range("I134").Value = 5330.778
$ws.Range("J134").Value = 3416.6667
$ws.Range("K134").Value = 15992.334
$ws.Range("L134").Value = 10250.0001
$ws.Range("M134").Value = -10922.334
$ws.Range("N134").Value = -20390.0001
$ws.Range("H136").Value = 1700.3636
$ws.Range("I136").Value = 1501
$ws.Range("K136").Value = 4503
$ws.Range("M136").Value = 597
$ws.Range("H138").Value = 5240
$ws.Range("I138").Value = 5240
$ws.Range("K138").Value = 15720
$ws.Range("M138").Value = -10580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 100005150
$ws.Range("I70").Value = 4785.143
$ws.Range("K70").Value = 4785.143
$ws.Range("M70").Value = -4515.143
$ws.Range("H73").Value = 100005150
$ws.Range("I73").Value = 4785.143
$ws.Range("K73").Value = 4785.143
$ws.Range("M73").Value = -3849.143
$ws.Range("H102").Value = 3712.25
$ws.Range("I102").Value = 3323.889
$ws.Range("K102").Value = 3323.889
$ws.Range("M102").Value = -1701.889
$ws.Range("H122").Value = 3082.2856
$ws.Range("I122").Value = 2883.1667
$ws.Range("K122").Value = 8649.500100000001
$ws.Range("M122").Value = -6199.500100000001
$ws.Range("H132").Value = 30335.025
$ws.Range("I132").Value = 45741.707
$ws.Range("K132").Value = 137225.121
$ws.Range("M132").Value = -134695.121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 45000
$ws.Range("J47").Value = 45000
$ws.Range("L47").Value = 45000
$ws.Range("N47").Value = -46144
$ws.Range("H97").Value = 49995
$ws.Range("J97").Value = 49995
$ws.Range("L97").Value = 49995
$ws.Range("N97").Value = -51977
$ws.Range("H136").Value = 4596.4707
$ws.Range("I136").Value = 3364.6
$ws.Range("J136").Value = 8018.3335
$ws.Range("K136").Value = 10093.8
$ws.Range("L136").Value = 24055.0005
$ws.Range("M136").Value = -7543.799999999999
$ws.Range("N136").Value = -29155.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 64499.5
$ws.Range("J95").Value = 64499.5
$ws.Range("L95").Value = 64499.5
$ws.Range("N95").Value = -69991.5
$ws.Range("H132").Value = 37039416
$ws.Range("I132").Value = 200001230
$ws.Range("K132").Value = 600003690
$ws.Range("M132").Value = -600001160

